$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 1191.14116875
$schedule.Range("F2").Value = 26.25972594246032

$schedule.Range("E3").Value = 367.04655
$schedule.Range("F3").Value = 24.27556547619048

# --- Sheet: Detailed ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B7").Value = 69.42238
$detailed.Range("B8").Value = 60.46039
$detailed.Range("B9").Value = 84.7901

$detailed.Range("B10").Value = 87.15796
$detailed.Range("C10").Value = "historical"

$detailed.Range("B13").Value = 93.85476

$detailed.Range("B16").Value = 50.66923

$detailed.Range("B18").Value = 56.98

$detailed.Range("B19").Value = 56.59029

$detailed.Range("B20").Value = 47.14626

$detailed.Range("B22").Value = 46.76395

$detailed.Range("B23").Value = 48.22486

$detailed.Range("B24").Value = 40.96027

$detailed.Range("B25").Value = 36.07

$detailed.Range("B27").Value = 40.82473

$detailed.Range("B32").Value = 27.27351

$detailed.Range("B33").Value = 18.777

$detailed.Range("B34").Value = 17.76526

$detailed.Range("B35").Value = 7.72659

$detailed.Range("B38").Value = -3.13019

$detailed.Range("B39").Value = -3.03124

$detailed.Range("B44").Value = 22.01959

$detailed.Range("B45").Value = 64.8901

$detailed.Range("B46").Value = 57.09

$detailed.Range("B47").Value = 57.52342

$detailed.Range("B48").Value = 58.1896

$detailed.Range("B49").Value = 57.03885
